$wb = $excel.ActiveWorkbook

$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")

# Row 7's status changes from "Ready for handoff" to "Handback transform failed"
# across all sheets that reference it (Overview!B7/C7, zh-cn!C7, de-de!C7).
$overviewSheet = $wb.Worksheets.Item("Overview")
$overviewSheet.Range("B7").Value = "Handback transform failed"
$overviewSheet.Range("C7").Value = "Handback transform failed"
$zhSheet.Range("C7").Value = "Handback transform failed"
$deSheet.Range("C7").Value = "Handback transform failed"

# New "Error Detail" entries for row 7 on the zh-cn and de-de sheets.
$zhSheet.Range("K7").Value = "Handback file name: yq041egs.rcs is different with handoff file name: 398a085f-eaf0-4fe7-8fea-7733b3c80634.3b95c3322b1e1d8bcb1d8c53caa9e73b448ed0a3.zh-cn."
$deSheet.Range("K7").Value = "Handback file name: yq041egs.rcs is different with handoff file name: 398a085f-eaf0-4fe7-8fea-7733b3c80634.3b95c3322b1e1d8bcb1d8c53caa9e73b448ed0a3.de-de."
